$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.771.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.909.70"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.00%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.72"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.47%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.503"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.913.47"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.78"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.144"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.51"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.45%  "
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.399.80"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.776.50"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.73"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.916.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.680"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.80"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.87"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.19"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.16"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0864"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -2.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.65"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.76"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.98"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.77%  "
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.67"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.294"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("E45").Value = "  -3.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "373.03"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.667.90"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.82"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.32"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.106"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.90%  "
